# The document ends with a trailing empty paragraph (just a paragraph
# mark, no runs). Turn it into a bold paragraph announcing the Python
# program file name, matching the target OOXML exactly:
#
#   <w:p>
#     <w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>
#     <w:r><w:rPr><w:b/><w:bCs/></w:rPr>
#       <w:t>Nom du programme Python :  READ_CSV_FINAL.py </w:t>
#     </w:r>
#   </w:p>

$d = $word.ActiveDocument

$target = $d.Paragraphs.Last
$range = $target.Range

$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
              '<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +
              '<w:r><w:rPr><w:b/><w:bCs/></w:rPr>' +
              '<w:t xml:space="preserve">Nom du programme Python :  READ_CSV_FINAL.py </w:t></w:r>' +
              '</w:p>'

# InsertXML replaces the contents of the target range (here, the whole
# trailing empty paragraph) with the raw OOXML fragment, giving us exact
# control over both the paragraph-mark run properties (pPr/rPr) and the
# run's own run properties (r/rPr).
$range.InsertXML($newParaXml)
